# fix product UTF8 encodage
# Replace accented French product names with plain-ASCII transliterations.
# Order matters: new shared-string entries are appended in the order the
# cell values are (re)written, and the target workbook appends them in the
# order the rows appear in the sheet (row 5, 15, 20, 27, 29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value  = "Vetements de travailleur"   # was "Vêtements de travailleur"
$ws.Range("B15").Value = "Biere"                      # was "Bière"
$ws.Range("B20").Value = "Montre a gousset"            # was "Montre à gousset"
$ws.Range("B27").Value = "Ble"                         # was "Blé"
$ws.Range("B29").Value = "Viande sechee"               # was "Viande séchée"
